# leaderboard.xlsx update: "Cycle 4 trials + updated leaderboard"
#
# 1) Fixes capitalization of a couple of existing player-name cells
#    ("Not Logic" -> "NotLogic", "tmty" -> "tTty") and (re)applies the
#    centered/word-wrap "Player name" cell style used throughout column A.
# 2) Appends 9 new Cycle-4 trial rows (152-160) with Player / Cycle /
#    Video / Time / Character1-8 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows: player-name corrections ---------------------------
$ws.Cells.Item(144, 1).Value = "NotLogic"
$ws.Cells.Item(147, 1).Value = "NotLogic"
$ws.Cells.Item(149, 1).Value = "NotLogic"
$ws.Cells.Item(151, 1).Value = "tTty"
$ws.Cells.Item(151, 3).Value = "https://youtu.be/wcdNh5GPRgo?si=lYENtBPqVkZpbsZD"

# Re-apply the standard centered / wrap-text player-name format (matches
# the style already used on I2 and throughout column A) to the 3 cells
# whose text just changed.
$ws.Range("I2").Copy()
$ws.Range("A144").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("A147").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("A149").PasteSpecial(-4122)

# --- New rows: Cycle 4 trial entries -----------------------------------
# Row 152: player "Echidna"
$ws.Cells.Item(152,1).Value = "Echidna"
$ws.Cells.Item(152,2).Value = 2
$ws.Cells.Item(152,3).Value = "https://youtu.be/6mDkkMd6idE"
$ws.Cells.Item(152,4).Value = 0.1986111111111111
$ws.Cells.Item(152,5).Value = "Hutao"
$ws.Cells.Item(152,6).Value = "Yelan"
$ws.Cells.Item(152,7).Value = "Mona"
$ws.Cells.Item(152,8).Value = "Zhongli"
$ws.Cells.Item(152,9).Value = "Neuvillette"
$ws.Cells.Item(152,10).Value = "Furina"
$ws.Cells.Item(152,11).Value = "Baizhu"
$ws.Cells.Item(152,12).Value = "Dehya"

# Row 153: player "Minishooo"
$ws.Cells.Item(153,1).Value = "Minishooo"
$ws.Cells.Item(153,2).Value = 3
$ws.Cells.Item(153,3).Value = "https://youtu.be/8SykPtYM5mg"
$ws.Cells.Item(153,4).Value = 0.18263888888888888
$ws.Cells.Item(153,5).Value = "Ayaka"
$ws.Cells.Item(153,6).Value = "Mona"
$ws.Cells.Item(153,7).Value = "Diona"
$ws.Cells.Item(153,8).Value = "Kazuha"
$ws.Cells.Item(153,9).Value = "Navia"
$ws.Cells.Item(153,10).Value = "Furina"
$ws.Cells.Item(153,11).Value = "Yelan"
$ws.Cells.Item(153,12).Value = "Jean"

# Row 154: player "NotLogic"
$ws.Cells.Item(154,1).Value = "NotLogic"
$ws.Cells.Item(154,2).Value = 3
$ws.Cells.Item(154,3).Value = "https://youtu.be/gfOq1Wga8qU"
$ws.Cells.Item(154,4).Value = 0.13541666666666666
$ws.Cells.Item(154,5).Value = "Neuvillette"
$ws.Cells.Item(154,6).Value = "Furina"
$ws.Cells.Item(154,7).Value = "Jean"
$ws.Cells.Item(154,8).Value = "Dehya"
$ws.Cells.Item(154,9).Value = "Hutao"
$ws.Cells.Item(154,10).Value = "Yelan"
$ws.Cells.Item(154,11).Value = "Mona"
$ws.Cells.Item(154,12).Value = "Zhongli"

# Row 155: player "Staryy"
$ws.Cells.Item(155,1).Value = "Staryy"
$ws.Cells.Item(155,2).Value = 3
$ws.Cells.Item(155,3).Value = "https://youtu.be/Wkv7RYaS-Qs"
$ws.Cells.Item(155,4).Value = 0.12083333333333333
$ws.Cells.Item(155,5).Value = "Ayaka"
$ws.Cells.Item(155,6).Value = "Mona"
$ws.Cells.Item(155,7).Value = "Shenhe"
$ws.Cells.Item(155,8).Value = "Sucrose"
$ws.Cells.Item(155,9).Value = "Yoimiya"
$ws.Cells.Item(155,10).Value = "Xingqiu"
$ws.Cells.Item(155,11).Value = "Yelan"
$ws.Cells.Item(155,12).Value = "YunJin"

# Row 156: player "Kyeburr"
$ws.Cells.Item(156,1).Value = "Kyeburr"
$ws.Cells.Item(156,2).Value = 3
$ws.Cells.Item(156,3).Value = "https://www.youtube.com/watch?v=z7Yig4M1Lfs"
$ws.Cells.Item(156,4).Value = 0.15347222222222223
$ws.Cells.Item(156,5).Value = "Lisa"
$ws.Cells.Item(156,6).Value = "Kazuha"
$ws.Cells.Item(156,7).Value = "Nahida"
$ws.Cells.Item(156,8).Value = "Kuki"
$ws.Cells.Item(156,9).Value = "Navia"
$ws.Cells.Item(156,10).Value = "TravelerGeo"
$ws.Cells.Item(156,11).Value = "Yae"
$ws.Cells.Item(156,12).Value = "Fischl"

# Row 157: player "Tmty"
$ws.Cells.Item(157,1).Value = "Tmty"
$ws.Cells.Item(157,2).Value = 3
$ws.Cells.Item(157,3).Value = "https://youtu.be/f2hYwmHafPo"
$ws.Cells.Item(157,4).Value = 0.19166666666666668
$ws.Cells.Item(157,5).Value = "Nilou"
$ws.Cells.Item(157,6).Value = "Yaoyao"
$ws.Cells.Item(157,7).Value = "Collei"
$ws.Cells.Item(157,8).Value = "Barbara"
$ws.Cells.Item(157,9).Value = "Yoimiya"
$ws.Cells.Item(157,10).Value = "Furina"
$ws.Cells.Item(157,11).Value = "Bennett"
$ws.Cells.Item(157,12).Value = "Diona"

# Row 158: player "Ghosted"
$ws.Cells.Item(158,1).Value = "Ghosted"
$ws.Cells.Item(158,2).Value = 3
$ws.Cells.Item(158,3).Value = "https://youtu.be/DJ2uIpoVHQw"
$ws.Cells.Item(158,4).Value = 0.22361111111111112
$ws.Cells.Item(158,5).Value = "Neuvillette"
$ws.Cells.Item(158,6).Value = "Furina"
$ws.Cells.Item(158,7).Value = "Jean"
$ws.Cells.Item(158,8).Value = "Beidou"
$ws.Cells.Item(158,9).Value = "Allhaitham"
$ws.Cells.Item(158,10).Value = "Fischl"
$ws.Cells.Item(158,11).Value = "Yae"
$ws.Cells.Item(158,12).Value = "Zhongli"

# Row 159: player "Tmty"
$ws.Cells.Item(159,1).Value = "Tmty"
$ws.Cells.Item(159,2).Value = 3
$ws.Cells.Item(159,3).Value = "https://youtu.be/OVxZHxHTuPA?si=Q1NztQm6t4ao1U1z"
$ws.Cells.Item(159,4).Value = 0.15625
$ws.Cells.Item(159,5).Value = "Ganyu"
$ws.Cells.Item(159,6).Value = "Mona"
$ws.Cells.Item(159,7).Value = "Rosaria"
$ws.Cells.Item(159,8).Value = "Kazuha"
$ws.Cells.Item(159,9).Value = "Yoimiya"
$ws.Cells.Item(159,10).Value = "Furina"
$ws.Cells.Item(159,11).Value = "Bennett"
$ws.Cells.Item(159,12).Value = "Diona"

# Row 160: player "Yuenn"
$ws.Cells.Item(160,1).Value = "Yuenn"
$ws.Cells.Item(160,2).Value = 3
$ws.Cells.Item(160,3).Value = "https://youtu.be/NC2e3XpPf7U?si=WQr_1tWRjXtftghU"
$ws.Cells.Item(160,4).Value = 0.19652777777777777
$ws.Cells.Item(160,5).Value = "Ayaka"
$ws.Cells.Item(160,6).Value = "Diona"
$ws.Cells.Item(160,7).Value = "Kazuha"
$ws.Cells.Item(160,8).Value = "Mona"
$ws.Cells.Item(160,9).Value = "Yoimiya"
$ws.Cells.Item(160,10).Value = "Xingqiu"
$ws.Cells.Item(160,11).Value = "Yelan"
$ws.Cells.Item(160,12).Value = "YunJin"

# Give the new rows' Time column (D) the same h:mm time-of-day format
# used by the rest of the table.
$ws.Range("D144").Copy()
$ws.Range("D152:D160").PasteSpecial(-4122)

# Update the window selection to match the bottom of the (now longer)
# table, mirroring where the author left the cursor after entering data.
$ws.Range("H166").Select()
